$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.325597167015076
$ws.Range("B1").Value = 1.706166982650757
$ws.Range("C1").Value = 6.685916900634766
$ws.Range("D1").Value = 1.610659122467041
$ws.Range("E1").Value = 0.9554146528244019
